# Apply targeted corrections to stock-qty / amount figures (and their
# roll-up Sub Total / Grand Total rows) on the single worksheet of the
# Companywise Stock Report, per the reconciled source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - 3M-SB FOOTLOCK MOP
$ws.Range("F5").Value = 13
$ws.Range("G5").Value = 4898.53

# Row 10 - Sub Total:
$ws.Range("B10").Value = 25818.35

# Row 30 - ZOFF Chaat Masala 100 GMS
$ws.Range("F30").Value = 12
$ws.Range("G30").Value = 423.48

# Row 34 - ZOFF Jeera Powder 100 GMS
$ws.Range("F34").Value = 67
$ws.Range("G34").Value = 2988.2

# Row 38 - Sub Total:
$ws.Range("B38").Value = 12094.12

# Row 77 - BHA-Whisper Ultra Clean XL Wings 7s pack
$ws.Range("F77").Value = 289
$ws.Range("G77").Value = 18409.3

# Row 96 - Sub Total:
$ws.Range("B96").Value = 160851.34

# Row 177 - GHP-Glamic Cistem Toilet Cleaner Cube
$ws.Range("F177").Value = 208
$ws.Range("G177").Value = 13478.4

# Row 181 - GHP-Glamic Disinfactant surface cleaner 500ML
$ws.Range("F181").Value = 16
$ws.Range("G181").Value = 743.84

# Row 184 - Sub Total:
$ws.Range("B184").Value = 27791.64

# Row 220 - HIM-GENTLE BABY SOAP 75G
$ws.Range("B220").Value = 48706
$ws.Range("E220").Value = 39.8
$ws.Range("F220").Value = -144
$ws.Range("G220").Value = -4795.2

# Row 221 - HIM-GENTLE BABY SOAP 75G
$ws.Range("B221").Value = 64973
$ws.Range("E221").Value = 35.4
$ws.Range("F221").Value = 0
$ws.Range("G221").Value = 0

# Row 231 - HIM-LIP BALM (12S BLISTER PACK ) 10G
$ws.Range("F231").Value = 47
$ws.Range("G231").Value = 947.52

# Row 244 - Sub Total:
$ws.Range("B244").Value = 29153.85

# Row 247 - HUL-Bru Inst Poly 50g
$ws.Range("B247").Value = 63565
$ws.Range("E247").Value = 109.19
$ws.Range("F247").Value = 60
$ws.Range("G247").Value = 6162.6

# Row 248 - HUL-Bru Inst Poly 50g
$ws.Range("B248").Value = 61610
$ws.Range("E248").Value = 122.71
$ws.Range("F248").Value = -58
$ws.Range("G248").Value = -5957.18

# Row 325 - JYOTHY-T-shine floor splst 500ml
$ws.Range("F325").Value = 24
$ws.Range("G325").Value = 1687.68

# Row 326 - JYOTHY-T-shine Toilet Liq 500ml lav
$ws.Range("F326").Value = 17
$ws.Range("G326").Value = 1213.12

# Row 334 - JYT - Ujala IDD Front load Detergent Liquid 2lt
$ws.Range("F334").Value = 20
$ws.Range("G334").Value = 4217.8

# Row 335 - JYT - Ujala IDD Top load Detergent Liquid 2lt
$ws.Range("F335").Value = 13
$ws.Range("G335").Value = 2741.57

# Row 336 - Sub Total:
$ws.Range("B336").Value = 156542.53

# Row 354 - KUS-Floor Wiper
$ws.Range("B354").Value = 58047
$ws.Range("D354").Value = 105.54
$ws.Range("E354").Value = 126.1
$ws.Range("F354").Value = 39
$ws.Range("G354").Value = 4116.06

# Row 355 - KUS-Floor Wiper
$ws.Range("B355").Value = 47097
$ws.Range("D355").Value = 112.28
$ws.Range("E355").Value = 134.16
$ws.Range("F355").Value = 15
$ws.Range("G355").Value = 1684.2

# Row 360 - KUS-Rope - 15Mtr (Plastic)
$ws.Range("F360").Value = 29
$ws.Range("G360").Value = 1079.09

# Row 362 - Sub Total:
$ws.Range("B362").Value = 24701.13

# Row 366 - CHUK-Black Pepper Powder 50GM
$ws.Range("F366").Value = 185
$ws.Range("G366").Value = 9586.700000000001

# Row 374 - CHUK-Lal mirch kutti 100gm
$ws.Range("F374").Value = 135
$ws.Range("G374").Value = 4275.45

# Row 377 - CHUK-Saunf Barik 100gm
$ws.Range("F377").Value = 24
$ws.Range("G377").Value = 1473.84

# Row 378 - Sub Total:
$ws.Range("B378").Value = 22636.09

# Row 414 - CRE-Butter cremfills 100gm
$ws.Range("B414").Value = 53263
$ws.Range("E414").Value = 15.29
$ws.Range("F414").Value = -309
$ws.Range("G414").Value = -3958.29

# Row 415 - CRE-Butter cremfills 100gm
$ws.Range("B415").Value = 65066
$ws.Range("E415").Value = 13.61
$ws.Range("F415").Value = 90
$ws.Range("G415").Value = 1152.9

# Row 423 - CRE-Cremica Honey Oatmeal Cookies 50 +25 Gm
$ws.Range("B423").Value = 64927
$ws.Range("E423").Value = 17.26
$ws.Range("F423").Value = 106
$ws.Range("G423").Value = 1719.32

# Row 424 - CRE-Cremica Honey Oatmeal Cookies 50 +25 Gm
$ws.Range("B424").Value = 45718
$ws.Range("E424").Value = 19.38
$ws.Range("F424").Value = -294
$ws.Range("G424").Value = -4768.68

# Row 428 - CRE-Cremica Oatmeal Digestive 112.5 Gm
$ws.Range("B428").Value = 45709
$ws.Range("E428").Value = 15.69
$ws.Range("F428").Value = -300
$ws.Range("G428").Value = -3945

# Row 429 - CRE-Cremica Oatmeal Digestive 112.5 Gm
$ws.Range("B429").Value = 64925
$ws.Range("E429").Value = 13.97
$ws.Range("F429").Value = 111
$ws.Range("G429").Value = 1459.65

# Row 490 - PRI-B-50 VIMAL Copper Glass 300ML (2pc Set)
$ws.Range("B490").Value = 64810
$ws.Range("E490").Value = 291.22
$ws.Range("F490").Value = 4
$ws.Range("G490").Value = 1095.68

# Row 491 - PRI-B-50 VIMAL Copper Glass 300ML (2pc Set)
$ws.Range("B491").Value = 53319
$ws.Range("E491").Value = 310.64
$ws.Range("F491").Value = -6
$ws.Range("G491").Value = -1643.52

# Row 503 - RAJ-Handy Brush (S-22)
$ws.Range("F503").Value = 39
$ws.Range("G503").Value = 2480.79

# Row 508 - Sub Total:
$ws.Range("B508").Value = 11837.18

# Row 525 - Nibbles Savoury (Jeera) 150G
$ws.Range("F525").Value = 4
$ws.Range("G525").Value = 181.36

# Row 526 - Sub Total:
$ws.Range("B526").Value = 181.36

# Row 603 - TCP-Coriander powder 200gm
$ws.Range("F603").Value = 21
$ws.Range("G603").Value = 1005.06

# Row 609 - TCP-TATA SAMPANN CHILLI POWDER 200 gm
$ws.Range("F609").Value = 10
$ws.Range("G609").Value = 695.6

# Row 614 - Sub Total:
$ws.Range("B614").Value = 2377.84

# Row 631 - Shankys Tip Top Vermicilli 500 Gm
$ws.Range("F631").Value = 32
$ws.Range("G631").Value = 1375.68

# Row 636 - Tip Top Sooji 1 Kg
$ws.Range("F636").Value = 15
$ws.Range("G636").Value = 854.85

# Row 637 - Sub Total:
$ws.Range("B637").Value = 11477.02

# Row 653 - VVD Priyam Cold Pressed Groundnut Oil Pouch 1 Ltr
$ws.Range("F653").Value = 1193
$ws.Range("G653").Value = 194590.23

# Row 655 - VVD Pure Drop Cold Pressed Gingelly Oil Pouch 1000Ml
$ws.Range("F655").Value = 345
$ws.Range("G655").Value = 97590.14999999999

# Row 656 - VVD Pure Drop Cold Pressed Gingelly Oil Pouch 500Ml
$ws.Range("F656").Value = 299
$ws.Range("G656").Value = 43250.35

# Row 660 - Sub Total:
$ws.Range("B660").Value = 336278.78

# Row 679 - Sub Total:
$ws.Range("B679").Value = 3296963.76

# Row 680 - Grand Total:
$ws.Range("B680").Value = 3296963.76
